$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh each coin's Price (D) and 1h Volume change (E) with the latest scrape.
# Price cells are forced to text (leading apostrophe) so Excel doesn't
# reinterpret numeric-looking strings like "556.09" as numbers.

$ws.Range("D2").Value = "'63.763.16"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "'3.100.57"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'556.09"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'138.11"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'3.093.75"
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("D9").Value = "'0.489"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "'6.73"
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").Value = "'0.453"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "'35.41"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "'0.0000215"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "'3.600.12"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").Value = "'63.891.31"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "'3.110.43"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("D19").Value = "'503.12"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("D20").Value = "'6.64"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "'13.66"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "'0.702"
$ws.Range("E22").Value = "  +2.67%  "
$ws.Range("D23").Value = "'7.22"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").Value = "'12.34"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "'77.35"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").Value = "'2.77"
$ws.Range("E27").Value = "  +3.07%  "
$ws.Range("D28").Value = "'8.38"
$ws.Range("E28").Value = "  +4.24%  "
$ws.Range("D29").Value = "'2.04"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "'26.09"
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("D32").Value = "'2.54"
$ws.Range("E32").Value = "  -3.04%  "
$ws.Range("D33").Value = "'1.11"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").Value = "'541.67"
$ws.Range("E34").Value = "  -7.39%  "
$ws.Range("D35").Value = "'55.03"
$ws.Range("E35").Value = "  +5.74%  "
$ws.Range("D36").Value = "'5.87"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("D37").Value = "'5.19"
$ws.Range("E37").Value = "  -3.75%  "
$ws.Range("D38").Value = "'0.0413"
$ws.Range("E38").Value = "  +3.61%  "
$ws.Range("D39").Value = "'0.0798"
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("D40").Value = "'3.064.69"
$ws.Range("E40").Value = "  +4.07%  "
$ws.Range("D41").Value = "'0.118"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").Value = "'8.08"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").Value = "'2.58"
$ws.Range("E43").Value = "  -11.14%  "
$ws.Range("D44").Value = "'0.254"
$ws.Range("E44").Value = "  +4.75%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'2.09"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "'121.05"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("D48").Value = "'24.17"
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("D49").Value = "'0.106"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").Value = "'0.0₃0500"
$ws.Range("E50").Value = "  -4.93%  "
$ws.Range("D51").Value = "'2.02"
$ws.Range("E51").Value = "  -2.07%  "
